$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row - new columns I and J
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Copy the header style from an existing header cell (e.g. H1) so the new
# headers match the bold/bordered/centered formatting used by the rest of row 1.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Data rows for column I (I0) and J (IF)
$ws.Range("I2").Value = 1
$ws.Range("I3").Value = 6
$ws.Range("I4").Value = 7
$ws.Range("I5").Value = 7
$ws.Range("I6").Value = 6
$ws.Range("I7").Value = 7
$ws.Range("I8").Value = 6
$ws.Range("I9").Value = 8
$ws.Range("I10").Value = 7

$ws.Range("J2").Value = 1
$ws.Range("J3").Value = 6
$ws.Range("J4").Value = 8
$ws.Range("J5").Value = 7
$ws.Range("J6").Value = 7
$ws.Range("J7").Value = 8
$ws.Range("J8").Value = 7
$ws.Range("J9").Value = 8
$ws.Range("J10").Value = 7
